{"js": "// Update the date line above the table and the 25 multiplication\n// equations inside the table to the new values from the commit.\n//\n// Each piece of text is replaced via a *scoped* search (paragraph body\n// for the date, individual table-cell body for each equation) rather\n// than a whole-document body.search(), so only the specific run that\n// holds the text is touched (formatting such as rFonts/sz is kept\n// because we use insertText(..., \"Replace\") on the found range, which\n// reuses the existing run's formatting).\n\n// --- 1. Date heading: \"2025-07-19 Saturday\" -> \"2025-07-20 Sunday\" ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\nconst dateResults = dateParagraph.search(\"2025-07-19 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length === 0) {\n  throw new Error(\"Date text not found\");\n}\ndateResults.items[0].insertText(\"2025-07-20 Sunday\", \"Replace\");\nawait context.sync();\n\n// --- 2. Table of equations: 5 rows (0, 4, 9, 14, 19) x 5 columns ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// old -> new for each (row, col); row indices are the table's actual\n// (0-based) row indices, which include the blank spacer rows.\nconst cellEdits = [\n  { row: 0, col: 0, oldText: \"399\u00d72=798\", newText: \"505\u00d73=1515\" },\n  { row: 0, col: 1, oldText: \"888\u00d73=2664\", newText: \"552\u00d78=4416\" },\n  { row: 0, col: 2, oldText: \"667\u00d76=4002\", newText: \"286\u00d72=572\" },\n  { row: 0, col: 3, oldText: \"957\u00d75=4785\", newText: \"757\u00d79=6813\" },\n  { row: 0, col: 4, oldText: \"775\u00d76=4650\", newText: \"795\u00d73=2385\" },\n\n  { row: 4, col: 0, oldText: \"555\u00d72=1110\", newText: \"827\u00d78=6616\" },\n  { row: 4, col: 1, oldText: \"444\u00d79=3996\", newText: \"255\u00d72=510\" },\n  { row: 4, col: 2, oldText: \"162\u00d74=648\", newText: \"399\u00d78=3192\" },\n  { row: 4, col: 3, oldText: \"972\u00d72=1944\", newText: \"223\u00d76=1338\" },\n  { row: 4, col: 4, oldText: \"469\u00d79=4221\", newText: \"621\u00d72=1242\" },\n\n  { row: 9, col: 0, oldText: \"312\u00d79=2808\", newText: \"209\u00d74=836\" },\n  { row: 9, col: 1, oldText: \"870\u00d74=3480\", newText: \"527\u00d72=1054\" },\n  { row: 9, col: 2, oldText: \"816\u00d72=1632\", newText: \"192\u00d75=960\" },\n  { row: 9, col: 3, oldText: \"646\u00d75=3230\", newText: \"699\u00d78=5592\" },\n  { row: 9, col: 4, oldText: \"156\u00d77=1092\", newText: \"510\u00d78=4080\" },\n\n  { row: 14, col: 0, oldText: \"264\u00d79=2376\", newText: \"764\u00d75=3820\" },\n  { row: 14, col: 1, oldText: \"786\u00d73=2358\", newText: \"414\u00d73=1242\" },\n  { row: 14, col: 2, oldText: \"592\u00d78=4736\", newText: \"528\u00d77=3696\" },\n  { row: 14, col: 3, oldText: \"781\u00d77=5467\", newText: \"156\u00d77=1092\" },\n  { row: 14, col: 4, oldText: \"542\u00d79=4878\", newText: \"623\u00d74=2492\" },\n\n  { row: 19, col: 0, oldText: \"809\u00d78=6472\", newText: \"197\u00d78=1576\" },\n  { row: 19, col: 1, oldText: \"905\u00d72=1810\", newText: \"367\u00d78=2936\" },\n  { row: 19, col: 2, oldText: \"268\u00d75=1340\", newText: \"453\u00d79=4077\" },\n  { row: 19, col: 3, oldText: \"761\u00d78=6088\", newText: \"120\u00d76=720\" },\n  { row: 19, col: 4, oldText: \"565\u00d77=3955\", newText: \"257\u00d76=1542\" },\n];\n\nfor (const { row, col, oldText, newText } of cellEdits) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found in cell (${row}, ${col}): ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line above the table and the 25 multiplication\n# equations inside the table to the new values from the commit.\n#\n# Each replacement is scoped to the smallest relevant Range (the date\n# paragraph, or the specific table cell) rather than the whole\n# $d.Content story, so only the run holding that text is touched.\n# Find/Replace reuses the existing run's formatting (rFonts/sz), so no\n# manual formatting work is needed.\n\n$d = $word.ActiveDocument\n\n# --- 1. Date heading: \"2025-07-19 Saturday\" -> \"2025-07-20 Sunday\" ---\n$dateRange = $d.Paragraphs(1).Range\n$find = $dateRange.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2025-07-19 Saturday\"\n$find.Replacement.Text = \"2025-07-20 Sunday\"\n$find.Execute(\"2025-07-19 Saturday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2025-07-20 Sunday\", 2)\n\n# --- 2. Table of equations: 5 rows (1, 5, 10, 15, 20) x 5 columns ---\n# (Word COM table/row/column indices are 1-based.)\n$tbl = $d.Tables(1)\n\n$cellEdits = @(\n    @(1, 1, \"399\u00d72=798\", \"505\u00d73=1515\"),\n    @(1, 2, \"888\u00d73=2664\", \"552\u00d78=4416\"),\n    @(1, 3, \"667\u00d76=4002\", \"286\u00d72=572\"),\n    @(1, 4, \"957\u00d75=4785\", \"757\u00d79=6813\"),\n    @(1, 5, \"775\u00d76=4650\", \"795\u00d73=2385\"),\n\n    @(5, 1, \"555\u00d72=1110\", \"827\u00d78=6616\"),\n    @(5, 2, \"444\u00d79=3996\", \"255\u00d72=510\"),\n    @(5, 3, \"162\u00d74=648\", \"399\u00d78=3192\"),\n    @(5, 4, \"972\u00d72=1944\", \"223\u00d76=1338\"),\n    @(5, 5, \"469\u00d79=4221\", \"621\u00d72=1242\"),\n\n    @(10, 1, \"312\u00d79=2808\", \"209\u00d74=836\"),\n    @(10, 2, \"870\u00d74=3480\", \"527\u00d72=1054\"),\n    @(10, 3, \"816\u00d72=1632\", \"192\u00d75=960\"),\n    @(10, 4, \"646\u00d75=3230\", \"699\u00d78=5592\"),\n    @(10, 5, \"156\u00d77=1092\", \"510\u00d78=4080\"),\n\n    @(15, 1, \"264\u00d79=2376\", \"764\u00d75=3820\"),\n    @(15, 2, \"786\u00d73=2358\", \"414\u00d73=1242\"),\n    @(15, 3, \"592\u00d78=4736\", \"528\u00d77=3696\"),\n    @(15, 4, \"781\u00d77=5467\", \"156\u00d77=1092\"),\n    @(15, 5, \"542\u00d79=4878\", \"623\u00d74=2492\"),\n\n    @(20, 1, \"809\u00d78=6472\", \"197\u00d78=1576\"),\n    @(20, 2, \"905\u00d72=1810\", \"367\u00d78=2936\"),\n    @(20, 3, \"268\u00d75=1340\", \"453\u00d79=4077\"),\n    @(20, 4, \"761\u00d78=6088\", \"120\u00d76=720\"),\n    @(20, 5, \"565\u00d77=3955\", \"257\u00d76=1542\")\n)\n\nforeach ($edit in $cellEdits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $oldText = $edit[2]\n    $newText = $edit[3]\n\n    $cellRange = $tbl.Cell($row, $col).Range\n    $find = $cellRange.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
